$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 120279.1
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 122730.71
$ws.Range("K17").Value = 450
$ws.Range("L17").Value = 368192.13
$ws.Range("M17").Value = -282
$ws.Range("N17").Value = -368528.13
$ws.Range("H113").Value = 4178.4
$ws.Range("I113").Value = 3124.5
$ws.Range("J113").Value = 4561.636
$ws.Range("K113").Value = 3124.5
$ws.Range("L113").Value = 4561.636
$ws.Range("M113").Value = 129.5
$ws.Range("N113").Value = -11069.636
$ws.Range("H141").Value = 2976.7058
$ws.Range("I141").Value = 2906.1875
$ws.Range("K141").Value = 8718.5625
$ws.Range("M141").Value = -3538.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 25000
$ws.Range("J23").Value = 25000
$ws.Range("L23").Value = 25000
$ws.Range("N23").Value = -25518
$ws.Range("H32").Value = 3347.516
$ws.Range("I32").Value = 2893.6785
$ws.Range("J32").Value = 7583.3335
$ws.Range("K32").Value = 2893.6785
$ws.Range("L32").Value = 7583.3335
$ws.Range("M32").Value = -2606.6785
$ws.Range("N32").Value = -8157.3335
$ws.Range("H61").Value = 2146.6667
$ws.Range("I61").Value = 1720
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1720
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1508
$ws.Range("N61").Value = -3424
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H63").Value = 3736.5789
$ws.Range("I63").Value = 3398.75
$ws.Range("J63").Value = 3826.6667
$ws.Range("K63").Value = 3398.75
$ws.Range("L63").Value = 3826.6667
$ws.Range("M63").Value = -2712.75
$ws.Range("N63").Value = -5198.6667
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H66").Value = 3736.5789
$ws.Range("I66").Value = 3398.75
$ws.Range("J66").Value = 3826.6667
$ws.Range("K66").Value = 16993.75
$ws.Range("L66").Value = 19133.3335
$ws.Range("M66").Value = -13561.75
$ws.Range("N66").Value = -25997.3335
$ws.Range("H80").Value = 26666.666
$ws.Range("J80").Value = 26666.666
$ws.Range("L80").Value = 26666.666
$ws.Range("N80").Value = -28662.666
$ws.Range("H83").Value = 26666.666
$ws.Range("J83").Value = 26666.666
$ws.Range("L83").Value = 79999.99800000001
$ws.Range("N83").Value = -89983.99800000001
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2166.6667
$ws.Range("K102").Value = 2166.6667
$ws.Range("M102").Value = -544.6667000000002
$ws.Range("H136").Value = 2146.6667
$ws.Range("I136").Value = 1720
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5160
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2610
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1700.5416
$ws.Range("I86").Value = 1576.8235
$ws.Range("J86").Value = 2001
$ws.Range("K86").Value = 1576.8235
$ws.Range("L86").Value = 2001
$ws.Range("M86").Value = -453.8235
$ws.Range("N86").Value = -4247
$ws.Range("H89").Value = 1700.5416
$ws.Range("I89").Value = 1576.8235
$ws.Range("J89").Value = 2001
$ws.Range("K89").Value = 7884.1175
$ws.Range("L89").Value = 10005
$ws.Range("M89").Value = -2268.1175
$ws.Range("N89").Value = -21237
$ws.Range("H94").Value = 815.26666
$ws.Range("I94").Value = 753.5454999999999
$ws.Range("K94").Value = 753.5454999999999
$ws.Range("M94").Value = -302.5454999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3033.3333
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 4060
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 4060
$ws.Range("M62").Value = -1676
$ws.Range("N62").Value = -5308
$ws.Range("H65").Value = 3033.3333
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 4060
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 20300
$ws.Range("M65").Value = -8380
$ws.Range("N65").Value = -26540
$ws.Range("H99").Value = 4133.222
$ws.Range("I99").Value = 4499.5
$ws.Range("J99").Value = 4028.5715
$ws.Range("K99").Value = 4499.5
$ws.Range("L99").Value = 4028.5715
$ws.Range("M99").Value = -3001.5
$ws.Range("N99").Value = -7024.5715
$ws.Range("H105").Value = 1200
$ws.Range("I105").Value = 1200
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1200
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 547
$ws.Range("N105").ClearContents()
$ws.Range("H126").Value = 4133.222
$ws.Range("I126").Value = 4499.5
$ws.Range("J126").Value = 4028.5715
$ws.Range("K126").Value = 13498.5
$ws.Range("L126").Value = 12085.7145
$ws.Range("M126").Value = -11028.5
$ws.Range("N126").Value = -17025.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 60005.5
$ws.Range("I18").Value = 50005
$ws.Range("K18").Value = 50005
$ws.Range("M18").Value = -49712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 91729.09
$ws.Range("J22").Value = 902
$ws.Range("L22").Value = 902
$ws.Range("N22").Value = -1492
$ws.Range("H27").Value = 91729.09
$ws.Range("J27").Value = 902
$ws.Range("L27").Value = 902
$ws.Range("N27").Value = -1116
$ws.Range("H68").Value = 13230
$ws.Range("I68").Value = 16900
$ws.Range("J68").Value = 4666.6665
$ws.Range("K68").Value = 16900
$ws.Range("L68").Value = 4666.6665
$ws.Range("M68").Value = -16151
$ws.Range("N68").Value = -6164.6665
$ws.Range("H71").Value = 13230
$ws.Range("I71").Value = 16900
$ws.Range("J71").Value = 4666.6665
$ws.Range("K71").Value = 84500
$ws.Range("L71").Value = 23333.3325
$ws.Range("M71").Value = -80756
$ws.Range("N71").Value = -30821.3325
$ws.Range("H93").Value = 21319.285
$ws.Range("I93").Value = 1537.84
$ws.Range("J93").Value = 70772.89999999999
$ws.Range("K93").Value = 1537.84
$ws.Range("L93").Value = 70772.89999999999
$ws.Range("M93").Value = -289.8399999999999
$ws.Range("N93").Value = -73268.89999999999
$ws.Range("H122").Value = 3113.7144
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3113.7144
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 9341.143199999999
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -14241.1432
